# Update 'F' column (想去人数 / want-to-go count) values per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 1380
$ws.Cells.Item(3, 6).Value = 98
$ws.Cells.Item(5, 6).Value = 5721
$ws.Cells.Item(6, 6).Value = 450
$ws.Cells.Item(8, 6).Value = 3347
$ws.Cells.Item(9, 6).Value = 6552
$ws.Cells.Item(10, 6).Value = 195
$ws.Cells.Item(11, 6).Value = 1278
$ws.Cells.Item(12, 6).Value = 738
$ws.Cells.Item(13, 6).Value = 91
$ws.Cells.Item(15, 6).Value = 11
$ws.Cells.Item(16, 6).Value = 1108
$ws.Cells.Item(17, 6).Value = 23
$ws.Cells.Item(18, 6).Value = 90
$ws.Cells.Item(20, 6).Value = 159
$ws.Cells.Item(22, 6).Value = 934
$ws.Cells.Item(23, 6).Value = 10
$ws.Cells.Item(24, 6).Value = 29
$ws.Cells.Item(25, 6).Value = 7
$ws.Cells.Item(26, 6).Value = 95
$ws.Cells.Item(28, 6).Value = 1140
$ws.Cells.Item(30, 6).Value = 19
$ws.Cells.Item(31, 6).Value = 20
$ws.Cells.Item(32, 6).Value = 16
$ws.Cells.Item(33, 6).Value = 264
$ws.Cells.Item(35, 6).Value = 257
$ws.Cells.Item(36, 6).Value = 1160
$ws.Cells.Item(37, 6).Value = 50
$ws.Cells.Item(38, 6).Value = 82

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(11, 6).Value = 3
$ws.Cells.Item(20, 6).Value = 175
$ws.Cells.Item(21, 6).Value = 128
$ws.Cells.Item(24, 6).Value = 602
$ws.Cells.Item(28, 6).Value = 648
$ws.Cells.Item(29, 6).Value = 946
$ws.Cells.Item(30, 6).Value = 556
$ws.Cells.Item(32, 6).Value = 78
$ws.Cells.Item(34, 6).Value = 8
$ws.Cells.Item(36, 6).Value = 123

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 700
$ws.Cells.Item(5, 6).Value = 821
$ws.Cells.Item(6, 6).Value = 531
$ws.Cells.Item(8, 6).Value = 743

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 700
$ws.Cells.Item(4, 6).Value = 1380
$ws.Cells.Item(5, 6).Value = 821
$ws.Cells.Item(7, 6).Value = 98
$ws.Cells.Item(8, 6).Value = 531
$ws.Cells.Item(11, 6).Value = 509
$ws.Cells.Item(13, 6).Value = 5721
$ws.Cells.Item(14, 6).Value = 450
$ws.Cells.Item(16, 6).Value = 3347
$ws.Cells.Item(18, 6).Value = 6552
$ws.Cells.Item(19, 6).Value = 195
$ws.Cells.Item(20, 6).Value = 1278
$ws.Cells.Item(24, 6).Value = 738
$ws.Cells.Item(25, 6).Value = 91
$ws.Cells.Item(26, 6).Value = 1108
$ws.Cells.Item(27, 6).Value = 128
$ws.Cells.Item(28, 6).Value = 90
$ws.Cells.Item(29, 6).Value = 159
$ws.Cells.Item(31, 6).Value = 934
$ws.Cells.Item(32, 6).Value = 602
$ws.Cells.Item(33, 6).Value = 29
$ws.Cells.Item(34, 6).Value = 95
$ws.Cells.Item(35, 6).Value = 1140
$ws.Cells.Item(37, 6).Value = 20
$ws.Cells.Item(38, 6).Value = 946
$ws.Cells.Item(39, 6).Value = 16
$ws.Cells.Item(40, 6).Value = 556
$ws.Cells.Item(41, 6).Value = 264
$ws.Cells.Item(43, 6).Value = 78
$ws.Cells.Item(44, 6).Value = 257
$ws.Cells.Item(45, 6).Value = 8
$ws.Cells.Item(47, 6).Value = 123
$ws.Cells.Item(49, 6).Value = 82
